$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.176399999999997
$ws.Range("B9").Value = 8.662200000000006
$ws.Range("D11").Value = -8.399800000000004
$ws.Range("B18").Value = 4.875400000000003
$ws.Range("B20").Value = 5.607099999999998
